$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 387 (shifts existing rows 387+ down by one)
$ws.Rows.Item(387).Insert()
$ws.Range("B387").Value = 0.54940052574212994
$ws.Range("C387").Value = 0.76302814643841776

# Delete the two trailing junk rows (now at 403 and 404 after the insert)
$ws.Rows.Item(403).Delete()
$ws.Rows.Item(403).Delete()

# Scroll back to the top and select G6 (matches the final view state)
[void]$ws.Range("G6").Select()
